$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.324.45'
$ws.Range("E2").Value = '  +5.57%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.790.91'
$ws.Range("E3").Value = '  +3.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  +0.34%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.07'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  +0.22%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4904'
$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2698'
$ws.Range("E8").Value = '  +2.92%  '

$ws.Range("E9").Value = '  +1.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.790.95'
$ws.Range("E10").Value = '  +3.29%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '16.66'
$ws.Range("E11").Value = '  +3.20%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07046'
$ws.Range("E12").Value = '  +1.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.6327'
$ws.Range("E13").Value = '  +3.30%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.638'
$ws.Range("E14").Value = '  +2.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '80.36'
$ws.Range("E15").Value = '  +3.75%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '28.311.10'
$ws.Range("E16").Value = '  +6.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9994'
$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  +0.38%  '

$ws.Range("E19").Value = '  +0.78%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.09'
$ws.Range("E20").Value = '  +5.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.019.21'
$ws.Range("E21").Value = '  +3.24%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.560'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.794'
$ws.Range("E23").Value = '  +2.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.261'
$ws.Range("E24").Value = '  +2.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.58'
$ws.Range("E25").Value = '  +2.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.81'
$ws.Range("E26").Value = '  +2.80%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.867'
$ws.Range("E27").Value = '  +5.45%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '110.20'
$ws.Range("E28").Value = '  +3.46%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.394'
$ws.Range("E29").Value = '  +0.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.153'
$ws.Range("E30").Value = '  +4.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08280'
$ws.Range("E31").Value = '  +3.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.801'
$ws.Range("E32").Value = '  +3.23%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04932'
$ws.Range("E33").Value = '  +9.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.083'
$ws.Range("E34").Value = '  +7.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6619'
$ws.Range("E35").Value = '  +6.55%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.615'
$ws.Range("E36").Value = '  +0.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9481'
$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.599'
$ws.Range("E38").Value = '  +7.44%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.097'
$ws.Range("E39").Value = '  +2.58%  '

$ws.Range("B40").Value = 'FraxShare'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.981'
$ws.Range("E40").Value = '  +6.59%  '

$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.01564'
$ws.Range("E41").Value = '  +3.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9983'
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("E43").Value = '  +0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4020'
$ws.Range("E44").Value = '  +4.06%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.165'
$ws.Range("E45").Value = '  +3.68%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1215'
$ws.Range("E46").Value = '  +4.54%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05452'
$ws.Range("E47").Value = '  +1.33%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.073'
$ws.Range("E48").Value = '  +2.59%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.298'
$ws.Range("E49").Value = '  +4.60%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '30.80'
$ws.Range("E50").Value = '  +1.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '53.08'
$ws.Range("E51").Value = '  +2.43%  '
